$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11
$ws.Range("B3").Value = 10
$ws.Range("B4").Value = 17
$ws.Range("B5").Value = 9
